$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -0.02110517777220711
$ws.Range("H4").Value = -0.01932204413288176
$ws.Range("J4").Value = 0.005298518211624132
$ws.Range("C5").Value = -0.01038435439937417
$ws.Range("H5").Value = 0.001044555689782227
$ws.Range("J5").Value = -0.0240587194032336
$ws.Range("C6").Value = 0.005599445407977816
$ws.Range("H6").Value = -0.00433904628556185
$ws.Range("J6").Value = 0.007944272400270062
$ws.Range("C7").Value = 0.003709482100379283
$ws.Range("H7").Value = 0.005263369362534773
$ws.Range("J7").Value = 0.03091636404310185
$ws.Range("C8").Value = 0.1003134614205384
$ws.Range("H8").Value = 0.9999999260799969
$ws.Range("J8").Value = -0.005078108643618663
$ws.Range("C9").Value = 0.9556117012964678
$ws.Range("H9").Value = -0.02724838803393552
$ws.Range("J9").Value = -0.01420375172972474
$ws.Range("C10").Value = 0.004209853032394121
$ws.Range("H10").Value = 0.01121331462453258
$ws.Range("J10").Value = 0.0008725577868348304
$ws.Range("C11").Value = -0.002767619438704777
$ws.Range("H11").Value = 0.005150281358011254
$ws.Range("J11").Value = -0.0131357627083223
$ws.Range("C12").Value = 0.06583369319334771
$ws.Range("H12").Value = 0.02384368088974723
$ws.Range("J12").Value = 0.005299909018230406
$ws.Range("C13").Value = 0.1040556061782242
$ws.Range("H13").Value = -0.01319918164796726
$ws.Range("J13").Value = -0.003690461918066865
$ws.Range("C14").Value = -0.2009272996530919
$ws.Range("H14").Value = 0.02071691545267662
$ws.Range("J14").Value = 0.008379181160454727
$ws.Range("C15").Value = -0.01358549516741981
$ws.Range("H15").Value = 0.00203412401736496
$ws.Range("J15").Value = 0.002421837686559839
$ws.Range("C16").Value = -0.00573481827739273
$ws.Range("H16").Value = 0.02683685982547439
$ws.Range("J16").Value = 0.005196766652402004
$ws.Range("C17").Value = 0.008202958888118355
$ws.Range("H17").Value = 0.04119785166391406
$ws.Range("J17").Value = -0.007845941383957901
$ws.Range("C18").Value = 0.02768330241933209
$ws.Range("H18").Value = 0.01336923394276936
$ws.Range("J18").Value = 0.01234322590784038
$ws.Range("C19").Value = 0.01138995165559806
$ws.Range("H19").Value = 0.002057256082290243
$ws.Range("J19").Value = -0.01171673490626191
$ws.Range("C20").Value = 0.009783347815333913
$ws.Range("H20").Value = -0.0005815280872611235
$ws.Range("J20").Value = -0.02955310253406249
$ws.Range("C21").Value = 0.02397830707113228
$ws.Range("H21").Value = 0.02410917101236684
$ws.Range("J21").Value = -0.01716372760577314
$ws.Range("C22").Value = 0.01142114109684564
$ws.Range("H22").Value = -0.0007680211507208459
$ws.Range("J22").Value = 0.0228612576295299
$ws.Range("C23").Value = -0.008519305108772203
$ws.Range("H23").Value = -0.008496251379850054
$ws.Range("J23").Value = -0.001876748153259229
